$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.9999999947744862,
    0.9940415142150969,
    0.97271369209994196,
    0.96285143992385414,
    0.95339761710654047,
    0.93019289452123455,
    0.92895141626956645,
    0.92783911189247648,
    0.92993511611512081,
    0.93301248950780002,
    0.93186267768874798,
    0.92993116101800555,
    0.92250843467513199,
    0.91999236906642479,
    0.9189899063001048,
    0.91838851425527968,
    0.91859852446233625,
    0.91938419562445262,
    0.99013538970749038,
    0.95778106045693767,
    0.95638252542500957,
    0.95511800077061837,
    0.9656780274632446,
    0.95265626738472786,
    0.94619912119046146,
    0.92653597926722009,
    0.92168762336171817,
    0.90020884470422957,
    0.88493523035240484,
    0.87836434043783052,
    0.87071011687597,
    0.86903075947866237,
    0.86851074456960453
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
